$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "B2"; Value = 4034 },
    @{ Cell = "G2"; Value = 15355 },
    @{ Cell = "H2"; Value = 4007 },
    @{ Cell = "J2"; Value = 377 },
    @{ Cell = "K2"; Value = 793 },
    @{ Cell = "O2"; Value = 899 },
    @{ Cell = "P2"; Value = 2664 },
    @{ Cell = "R2"; Value = 8989 },
    @{ Cell = "S2"; Value = 12540 },
    @{ Cell = "U2"; Value = 6499 },
    @{ Cell = "V2"; Value = 4308 },
    @{ Cell = "AB2"; Value = 126240 },
    @{ Cell = "AC2"; Value = 11037 },
    @{ Cell = "AE2"; Value = 137929 },
    @{ Cell = "G3"; Value = 15215 },
    @{ Cell = "H3"; Value = 4035 },
    @{ Cell = "N3"; Value = 2328 },
    @{ Cell = "P3"; Value = 2480 },
    @{ Cell = "S3"; Value = 12832 },
    @{ Cell = "AB3"; Value = 128464 },
    @{ Cell = "AC3"; Value = 11227 },
    @{ Cell = "AE3"; Value = 140306 },
    @{ Cell = "G4"; Value = 15221 },
    @{ Cell = "H4"; Value = 3921 },
    @{ Cell = "O4"; Value = 842 },
    @{ Cell = "P4"; Value = 2494 },
    @{ Cell = "R4"; Value = 9192 },
    @{ Cell = "S4"; Value = 13055 },
    @{ Cell = "AB4"; Value = 131335 },
    @{ Cell = "AC4"; Value = 11505 },
    @{ Cell = "AD4"; Value = 634 },
    @{ Cell = "AE4"; Value = 143475 },
    @{ Cell = "G5"; Value = 15036 },
    @{ Cell = "H5"; Value = 3911 },
    @{ Cell = "O5"; Value = 898 },
    @{ Cell = "P5"; Value = 2464 },
    @{ Cell = "R5"; Value = 9512 },
    @{ Cell = "S5"; Value = 13204 },
    @{ Cell = "U5"; Value = 7315 },
    @{ Cell = "AB5"; Value = 133409 },
    @{ Cell = "AC5"; Value = 11747 },
    @{ Cell = "AE5"; Value = 145801 },
    @{ Cell = "B6"; Value = 4368 },
    @{ Cell = "G6"; Value = 15302 },
    @{ Cell = "H6"; Value = 4163 },
    @{ Cell = "K6"; Value = 839 },
    @{ Cell = "N6"; Value = 2195 },
    @{ Cell = "O6"; Value = 842 },
    @{ Cell = "P6"; Value = 2576 },
    @{ Cell = "S6"; Value = 13722 },
    @{ Cell = "U6"; Value = 7407 },
    @{ Cell = "V6"; Value = 5015 },
    @{ Cell = "AB6"; Value = 134979 },
    @{ Cell = "AC6"; Value = 12176 },
    @{ Cell = "AD6"; Value = 696 },
    @{ Cell = "AE6"; Value = 147852 },
    @{ Cell = "B7"; Value = 4451 },
    @{ Cell = "E7"; Value = 14108 },
    @{ Cell = "G7"; Value = 15892 },
    @{ Cell = "H7"; Value = 4243 },
    @{ Cell = "I7"; Value = 2141 },
    @{ Cell = "O7"; Value = 874 },
    @{ Cell = "P7"; Value = 2617 },
    @{ Cell = "Q7"; Value = 4255 },
    @{ Cell = "S7"; Value = 14237 },
    @{ Cell = "U7"; Value = 7586 },
    @{ Cell = "V7"; Value = 5374 },
    @{ Cell = "AB7"; Value = 139999 },
    @{ Cell = "AC7"; Value = 12672 },
    @{ Cell = "AD7"; Value = 748 },
    @{ Cell = "AE7"; Value = 153419 },
    @{ Cell = "B8"; Value = 4387 },
    @{ Cell = "G8"; Value = 15921 },
    @{ Cell = "H8"; Value = 4239 },
    @{ Cell = "I8"; Value = 2062 },
    @{ Cell = "J8"; Value = 367 },
    @{ Cell = "P8"; Value = 2705 },
    @{ Cell = "R8"; Value = 9725 },
    @{ Cell = "S8"; Value = 14230 },
    @{ Cell = "U8"; Value = 7742 },
    @{ Cell = "AB8"; Value = 141285 },
    @{ Cell = "AC8"; Value = 12809 },
    @{ Cell = "AE8"; Value = 154766 },
    @{ Cell = "B9"; Value = 4318 },
    @{ Cell = "G9"; Value = 15334 },
    @{ Cell = "H9"; Value = 4323 },
    @{ Cell = "J9"; Value = 290 },
    @{ Cell = "K9"; Value = 801 },
    @{ Cell = "N9"; Value = 2413 },
    @{ Cell = "O9"; Value = 827 },
    @{ Cell = "P9"; Value = 2640 },
    @{ Cell = "Q9"; Value = 4153 },
    @{ Cell = "R9"; Value = 8338 },
    @{ Cell = "S9"; Value = 13861 },
    @{ Cell = "U9"; Value = 6362 },
    @{ Cell = "AB9"; Value = 132889 },
    @{ Cell = "AC9"; Value = 11972 },
    @{ Cell = "AE9"; Value = 145498 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
